# Auto-generated edit script: update cryptos price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.764.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.65%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.648.15'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.19%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.27%  '

# Row 6
$ws.Range("E6").Value = '  +4.61%  '

# Row 7
$ws.Range("E7").Value = '  -0.08%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.23'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.27%  '

# Row 9
$ws.Range("E9").Value = '  +0.71%  '

# Row 10
$ws.Range("E10").Value = '  +0.31%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0891'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.42%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.879.74'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.30%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.640.77'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.64%  '

# Row 14
$ws.Range("E14").Value = '  -0.51%  '

# Row 15
$ws.Range("E15").Value = '  -0.36%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.59%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.722.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.49%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.04'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.96%  '

# Row 19
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0726'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.11%  '

# Row 20
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.00%  '

# Row 22
$ws.Range("E22").Value = '  -0.49%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.11%  '

# Row 24
$ws.Range("E24").Value = '  -3.64%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.98'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.88%  '

# Row 26
$ws.Range("E26").Value = '  -0.98%  '

# Row 27
$ws.Range("E27").Value = '  +1.35%  '

# Row 28
$ws.Range("E28").Value = '  -0.55%  '

# Row 29
$ws.Range("E29").Value = '  -0.13%  '

# Row 30
$ws.Range("E30").Value = '  +0.17%  '

# Row 31
$ws.Range("E31").Value = '  -1.61%  '

# Row 32
$ws.Range("E32").Value = '  +0.68%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.446.80'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.36%  '

# Row 34
$ws.Range("E34").Value = '  +1.31%  '

# Row 35
$ws.Range("E35").Value = '  +1.71%  '

# Row 36
$ws.Range("E36").Value = '  -1.09%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.573'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.69%  '

# Row 38
$ws.Range("E38").Value = '  -1.71%  '

# Row 39
$ws.Range("E39").Value = '  -0.29%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.898'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +14.24%  '

# Row 41
$ws.Range("E41").Value = '  -1.49%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.97%  '

# Row 43
$ws.Range("E43").Value = '  -0.04%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.47'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.49%  '

# Row 45
$ws.Range("E45").Value = '  -0.42%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.26'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.03%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.789.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.25%  '

# Row 48
$ws.Range("E48").Value = '  +3.62%  '

# Row 49
$ws.Range("E49").Value = '  -1.38%  '

# Row 50
$ws.Range("E50").Value = '  +1.70%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0997'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.47%  '
